# Adds 1x12 board prices to the BOM workbook (rows 76-79), plus the
# subtotal in Q95, and updates the sheet view's scroll position/selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Row 76 : H.0404.09300.SI style siding bearing (93" x65 qty) ---
$ws.Range("L76").Value = 96
$ws.Range("M76:M79").Formula = '=INT(L76/12) & " ft "'
$ws.Range("N76").Value = 65
$ws.Range("O76").Value = 7.2
$ws.Range("P76:P79").Formula = "=O76*N76"
$ws.Range("R76").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Row 77 ---
$ws.Range("L77").Value = 120
$ws.Range("N77").Value = 18
$ws.Range("O77").Value = 9
$ws.Range("R77").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Row 78 ---
$ws.Range("L78").Value = 144
$ws.Range("N78").Value = 3
$ws.Range("O78").Value = 11.4
$ws.Range("R78").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Row 79 ---
$ws.Range("L79").Value = 168
$ws.Range("N79").Value = 3
$ws.Range("O79").Value = 13.3
$ws.Range("R79").Value = "https://collinssawmillcompany.wordpress.com/price-list/"

# --- Subtotal for the new 1x12 rows ---
$ws.Range("Q95").Formula = "=SUM(P76:P79)"

# --- Update sheet view: scroll position and active selection ---
$app = $wb.Application
$win = $app.ActiveWindow
$win.ScrollRow = 51
$win.ScrollColumn = 7
$ws.Range("R76:R79").Select()
